$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.327376365661621
$ws.Range("B1").Value = 1.17270302772522
$ws.Range("C1").Value = 3.170324325561523
$ws.Range("D1").Value = 3.025180578231812
$ws.Range("E1").Value = 0.9158830642700195
